{"js": "// Replace the date line and the multiplication problems per the diff mapping.\n// Each old string in the document is unique, so a targeted search + Replace\n// for each one is sufficient (keeps the original run/paragraph formatting).\nconst replacements = [\n  [\"2024-10-29 Tuesday\", \"2024-10-30 Wednesday\"],\n  [\"174\u00d73=\", \"184\u00d79=\"],\n  [\"142\u00d75=\", \"148\u00d74=\"],\n  [\"586\u00d73=\", \"531\u00d79=\"],\n  [\"820\u00d76=\", \"592\u00d74=\"],\n  [\"420\u00d76=\", \"288\u00d75=\"],\n  [\"439\u00d74=\", \"720\u00d72=\"],\n  [\"731\u00d72=\", \"923\u00d75=\"],\n  [\"120\u00d73=\", \"275\u00d72=\"],\n  [\"309\u00d79=\", \"646\u00d72=\"],\n  [\"309\u00d74=\", \"227\u00d78=\"],\n  [\"417\u00d74=\", \"336\u00d72=\"],\n  [\"967\u00d73=\", \"275\u00d78=\"],\n  [\"411\u00d78=\", \"646\u00d79=\"],\n  [\"621\u00d73=\", \"241\u00d76=\"],\n  [\"136\u00d72=\", \"722\u00d79=\"],\n  [\"366\u00d75=\", \"390\u00d75=\"],\n  [\"357\u00d73=\", \"259\u00d72=\"],\n  [\"412\u00d73=\", \"920\u00d77=\"],\n  [\"689\u00d76=\", \"223\u00d79=\"],\n  [\"439\u00d78=\", \"600\u00d77=\"],\n  [\"320\u00d73=\", \"896\u00d75=\"],\n  [\"598\u00d75=\", \"465\u00d79=\"],\n  [\"784\u00d77=\", \"787\u00d78=\"],\n  [\"823\u00d76=\", \"742\u00d74=\"],\n  [\"333\u00d77=\", \"967\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the multiplication problems per the diff mapping.\n# Each old string in the document is unique, so a targeted Find/Replace for\n# each one is sufficient (keeps the original run/paragraph formatting).\n$d = $word.ActiveDocument\n\nfunction Replace-OneText($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $wdReplaceAll = 2\n    $wdFindContinue = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n\nReplace-OneText \"2024-10-29 Tuesday\" \"2024-10-30 Wednesday\"\nReplace-OneText \"174\u00d73=\" \"184\u00d79=\"\nReplace-OneText \"142\u00d75=\" \"148\u00d74=\"\nReplace-OneText \"586\u00d73=\" \"531\u00d79=\"\nReplace-OneText \"820\u00d76=\" \"592\u00d74=\"\nReplace-OneText \"420\u00d76=\" \"288\u00d75=\"\nReplace-OneText \"439\u00d74=\" \"720\u00d72=\"\nReplace-OneText \"731\u00d72=\" \"923\u00d75=\"\nReplace-OneText \"120\u00d73=\" \"275\u00d72=\"\nReplace-OneText \"309\u00d79=\" \"646\u00d72=\"\nReplace-OneText \"309\u00d74=\" \"227\u00d78=\"\nReplace-OneText \"417\u00d74=\" \"336\u00d72=\"\nReplace-OneText \"967\u00d73=\" \"275\u00d78=\"\nReplace-OneText \"411\u00d78=\" \"646\u00d79=\"\nReplace-OneText \"621\u00d73=\" \"241\u00d76=\"\nReplace-OneText \"136\u00d72=\" \"722\u00d79=\"\nReplace-OneText \"366\u00d75=\" \"390\u00d75=\"\nReplace-OneText \"357\u00d73=\" \"259\u00d72=\"\nReplace-OneText \"412\u00d73=\" \"920\u00d77=\"\nReplace-OneText \"689\u00d76=\" \"223\u00d79=\"\nReplace-OneText \"439\u00d78=\" \"600\u00d77=\"\nReplace-OneText \"320\u00d73=\" \"896\u00d75=\"\nReplace-OneText \"598\u00d75=\" \"465\u00d79=\"\nReplace-OneText \"784\u00d77=\" \"787\u00d78=\"\nReplace-OneText \"823\u00d76=\" \"742\u00d74=\"\nReplace-OneText \"333\u00d77=\" \"967\u00d74=\"\n"}
